$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source sheet stores every Price/Volume(1h) figure as literal text (the
# upstream scraper writes raw strings like "63.584.71" / "0.120", preserving
# thousands-dot grouping and trailing zeros that a real Number type would drop).
# Excel auto-detects plain numeric-looking strings and would silently convert them
# (e.g. "0.120" -> 0.12), so for those cells we briefly force Text format, assign
# the literal value, then restore the default "Normal" cell style so no formatting
# footprint is left behind.

$ws.Range("D2").Value = "63.584.71"
$ws.Range("E2").Value = "  +1.00%  "
$ws.Range("D3").Value = "3.319.14"
$ws.Range("E3").Value = "  +5.57%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.32%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "3.318.46"
$ws.Range("E8").Value = "  +5.61%  "
$ws.Range("E9").Value = "  +0.82%  "
$ws.Range("E10").Value = "  +2.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.52"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.45%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.467"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.29%  "
$ws.Range("E13").Value = "  +1.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.69"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.90%  "
$ws.Range("D15").Value = "3.871.18"
$ws.Range("E15").Value = "  +5.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.120"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").Value = "3.323.11"
$ws.Range("E17").Value = "  +5.77%  "
$ws.Range("D18").Value = "63.678.02"
$ws.Range("E18").Value = "  +1.22%  "
$ws.Range("E19").Value = "  +3.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "479.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.733"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.71%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.37%  "
$ws.Range("E27").Value = "  +2.64%  "
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.62%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.97%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.92"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.90%  "
$ws.Range("E33").Value = "  +1.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.53"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.32%  "
$ws.Range("E35").Value = "  +3.83%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.05"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.78%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.40"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.15%  "
$ws.Range("D38").Value = "0.0₃0746"
$ws.Range("E38").Value = "  +7.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0398"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.96%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "434.12"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.55%  "
$ws.Range("D41").Value = "3.086.17"
$ws.Range("E41").Value = "  +5.75%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.121"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +8.73%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.76"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.87%  "
$ws.Range("B44").Value = "Cosmos"
$ws.Range("C44").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.33"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.91%  "
$ws.Range("E45").Value = "  +2.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.21"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "36.99"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +15.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.29"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.77%  "
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("E50").Value = "  +1.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.29"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.67%  "
